$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two "Order" rows (13 and 14) entirely - they shift everything
# below up and the now-unused "Order"/"asc"/"desc" shared strings get pruned.
$ws.Rows("13:14").Delete()

# Add the new row 12: View / Defensive / Free Cash Flow to Firm / fcffgrowth / ~gt~ / -99
$ws.Range("A12").Value = "View"
$ws.Range("B12").Value = "Defensive"
$ws.Range("C12").Value = "Free Cash Flow to Firm"
$ws.Range("D12").Value = "fcffgrowth"
$ws.Range("E12").Value = "~gt~"
$ws.Range("F12").Value = -99

# Update the active selection to match the target state
$ws.Range("D13").Select()
